# Applies the "Updated symbol list" data refresh to the crypto price sheet.
# Price cells (column D) are stored as text in the workbook, so we force the
# cell's number format to Text ("@") before assigning the value - otherwise
# Excel will auto-convert the numeric-looking string into a real number and
# mangle the exact text representation (trailing zeros, scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "242.03"

# Row 3 - OKB
Set-TextValue "D3" "21.87"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.383"

# Row 5 - Cronos
Set-TextValue "D5" "0.05691"

# Row 6 - GateToken
Set-TextValue "D6" "3.411"

# Row 7 - KuCoinToken
Set-TextValue "D7" "6.286"

# Row 8 - FTXToken
Set-TextValue "D8" "1.127"
$ws.Range("E8").Value = "7FTXTokenFTT"

# Row 9 - MXToken
Set-TextValue "D9" "0.8072"

# Row 10 - now "One" (rows 10-18 shifted up by one coin / re-ranked)
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01074"
$ws.Range("E10").Value = "9OneONEBestin24h"

# Row 11 - now "WazirX"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1421"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12 - now "MandalaExchangeToken"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07279"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13 - now "LiechtensteinCryptoassetsExchange"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03082"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14 - now "BitrueCoin"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03095"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15 - now "BitMartToken"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09352"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16 - now "MCDex"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.904"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17 - now "BitForexToken"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001584"
$ws.Range("E17").Value = "16BitForexTokenBF"

# Row 18 - now "CoinExToken"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04805"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19 - TigerCash
Set-TextValue "D19" "0.006298"

# Row 20 - BitKan
Set-TextValue "D20" "0.0009980"

# Row 21 - HotbitToken
Set-TextValue "D21" "0.004061"

# Row 22 - NitroEx
Set-TextValue "D22" "0.0001500"

# Row 24 - BTSEToken
Set-TextValue "D24" "2.153"

# Row 26 - ProBitToken
Set-TextValue "D26" "0.1299"

# Row 27 - UpBots
$ws.Range("E27").Value = "26UpBotsUBXT"

# Row 40 - IDEX
Set-TextValue "D40" "0.03813"

# Row 41 - KickToken
Set-TextValue "D41" "0.006665"

# Row 43 - CEJI
Set-TextValue "D43" "0.002832"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.006484"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005615"

# Row 47 - CoinbaseStockToken
Set-TextValue "D47" "0.3901"

# Row 48 - BOLO
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
